$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = '@'
$r.Value = '69.547.39'
$r.Style = 'Normal'
$ws.Range("E2").Value = '  +0.26%  '

$r = $ws.Range("D3")
$r.NumberFormat = '@'
$r.Value = '3.791.35'
$r.Style = 'Normal'
$ws.Range("E3").Value = '  +1.31%  '

$r = $ws.Range("D4")
$r.NumberFormat = '@'
$r.Value = '1.00'
$r.Style = 'Normal'
$ws.Range("E4").Value = '  -0.17%  '

$r = $ws.Range("D5")
$r.NumberFormat = '@'
$r.Value = '616.19'
$r.Style = 'Normal'
$ws.Range("E5").Value = '  +0.06%  '

$r = $ws.Range("D6")
$r.NumberFormat = '@'
$r.Value = '178.22'
$r.Style = 'Normal'
$ws.Range("E6").Value = '  +1.41%  '

$r = $ws.Range("D7")
$r.NumberFormat = '@'
$r.Value = '3.788.72'
$r.Style = 'Normal'
$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  -0.31%  '

$r = $ws.Range("D10")
$r.NumberFormat = '@'
$r.Value = '0.167'
$r.Style = 'Normal'
$ws.Range("E10").Value = '  -0.30%  '

$r = $ws.Range("D11")
$r.NumberFormat = '@'
$r.Value = '6.55'
$r.Style = 'Normal'
$ws.Range("E11").Value = '  +4.92%  '

$r = $ws.Range("D12")
$r.NumberFormat = '@'
$r.Value = '0.487'
$r.Style = 'Normal'
$ws.Range("E12").Value = '  +0.63%  '

$r = $ws.Range("D13")
$r.NumberFormat = '@'
$r.Value = '40.05'
$r.Style = 'Normal'
$ws.Range("E13").Value = '  -0.80%  '

$r = $ws.Range("D14")
$r.NumberFormat = '@'
$r.Value = '0.0000255'
$r.Style = 'Normal'
$ws.Range("E14").Value = '  -0.54%  '

$r = $ws.Range("D15")
$r.NumberFormat = '@'
$r.Value = '4.419.76'
$r.Style = 'Normal'
$ws.Range("E15").Value = '  +1.23%  '

$r = $ws.Range("D16")
$r.NumberFormat = '@'
$r.Value = '3.784.99'
$r.Style = 'Normal'
$ws.Range("E16").Value = '  +0.73%  '

$r = $ws.Range("D17")
$r.NumberFormat = '@'
$r.Value = '69.647.03'
$r.Style = 'Normal'
$ws.Range("E17").Value = '  +0.21%  '

$r = $ws.Range("D18")
$r.NumberFormat = '@'
$r.Value = '7.59'
$r.Style = 'Normal'
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("E19").Value = '  -3.06%  '

$r = $ws.Range("D20")
$r.NumberFormat = '@'
$r.Value = '511.68'
$r.Style = 'Normal'
$ws.Range("E20").Value = '  +1.43%  '

$r = $ws.Range("D21")
$r.NumberFormat = '@'
$r.Value = '16.48'
$r.Style = 'Normal'
$ws.Range("E21").Value = '  -0.59%  '

$r = $ws.Range("D22")
$r.NumberFormat = '@'
$r.Value = '9.44'
$r.Style = 'Normal'
$ws.Range("E22").Value = '  -0.25%  '

$r = $ws.Range("D23")
$r.NumberFormat = '@'
$r.Value = '0.737'
$r.Style = 'Normal'
$ws.Range("E23").Value = '  +2.50%  '

$r = $ws.Range("D24")
$r.NumberFormat = '@'
$r.Value = '2.51'
$r.Style = 'Normal'
$ws.Range("E24").Value = '  +0.71%  '

$r = $ws.Range("D25")
$r.NumberFormat = '@'
$r.Value = '86.57'
$r.Style = 'Normal'
$ws.Range("E25").Value = '  +0.13%  '

$r = $ws.Range("D26")
$r.NumberFormat = '@'
$r.Value = '12.94'
$r.Style = 'Normal'
$ws.Range("E26").Value = '  -0.97%  '

$r = $ws.Range("D27")
$r.NumberFormat = '@'
$r.Value = '0.0000137'
$r.Style = 'Normal'
$ws.Range("E27").Value = '  +0.47%  '

$r = $ws.Range("D28")
$r.NumberFormat = '@'
$r.Value = '10.62'
$r.Style = 'Normal'
$ws.Range("E28").Value = '  -3.78%  '

$ws.Range("E29").Value = '  +0.50%  '

$r = $ws.Range("D30")
$r.NumberFormat = '@'
$r.Value = '2.56'
$r.Style = 'Normal'
$ws.Range("E30").Value = '  +3.88%  '

$r = $ws.Range("D31")
$r.NumberFormat = '@'
$r.Value = '3.01'
$r.Style = 'Normal'
$ws.Range("E31").Value = '  +3.89%  '

$r = $ws.Range("D32")
$r.NumberFormat = '@'
$r.Value = '8.10'
$r.Style = 'Normal'
$ws.Range("E32").Value = '  +4.69%  '

$r = $ws.Range("D33")
$r.NumberFormat = '@'
$r.Value = '31.35'
$r.Style = 'Normal'
$ws.Range("E33").Value = '  +1.84%  '

$ws.Range("E34").Value = '  +1.48%  '

$ws.Range("E35").Value = '  -0.29%  '

$r = $ws.Range("D36")
$r.NumberFormat = '@'
$r.Value = '1.06'
$r.Style = 'Normal'
$ws.Range("E36").Value = '  +0.05%  '

$r = $ws.Range("D37")
$r.NumberFormat = '@'
$r.Value = '6.17'
$r.Style = 'Normal'
$ws.Range("E37").Value = '  +1.02%  '

$ws.Range("E38").Value = '  +7.55%  '

$ws.Range("E39").Value = '  +3.00%  '

$r = $ws.Range("D40")
$r.NumberFormat = '@'
$r.Value = '461.05'
$r.Style = 'Normal'
$ws.Range("E40").Value = '  +9.81%  '

$ws.Range("E41").Value = '  -0.83%  '

$ws.Range("E42").Value = '  -1.43%  '

$ws.Range("E43").Value = '  +7.65%  '

$r = $ws.Range("D44")
$r.NumberFormat = '@'
$r.Value = '44.39'
$r.Style = 'Normal'
$ws.Range("E44").Value = '  -1.50%  '

$r = $ws.Range("D45")
$r.NumberFormat = '@'
$r.Value = '8.61'
$r.Style = 'Normal'
$ws.Range("E45").Value = '  -0.20%  '

$r = $ws.Range("D46")
$r.NumberFormat = '@'
$r.Value = '2.965.05'
$r.Style = 'Normal'
$ws.Range("E46").Value = '  -1.60%  '

$r = $ws.Range("D47")
$r.NumberFormat = '@'
$r.Value = '0.0362'
$r.Style = 'Normal'
$ws.Range("E47").Value = '  +0.84%  '

$r = $ws.Range("D48")
$r.NumberFormat = '@'
$r.Value = '27.39'
$r.Style = 'Normal'
$ws.Range("E48").Value = '  +1.07%  '

$r = $ws.Range("D50")
$r.NumberFormat = '@'
$r.Value = '139.06'
$r.Style = 'Normal'
$ws.Range("E50").Value = '  +0.59%  '

$r = $ws.Range("D51")
$r.NumberFormat = '@'
$r.Value = '2.48'
$r.Style = 'Normal'
$ws.Range("E51").Value = '  +1.06%  '
